# Edit slide 5 ("STGraph - Operations"), Content Placeholder shape:
#   - Fix "Search Algorithm<:" -> "Search algorithm:"
#   - Simplify the interval-math line (drop sub/superscripts), move it up a level,
#     and duplicate it as a second identical bullet
#   - Rename/un-bold "GraphNode2GraphNode" -> "Join strategy: Nested-Loop;"
#   - Drop "Furthermore" / "Join: Nested-Join Loop;" / "Temporal traversal: ..."
#     bullets, replacing them with a single new bullet
#     "Each time a traversal goes through a virtual edge" one level in

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$mathLine = "[tAs, tAe) ∩ [tBs, tBe) ≠ ∅ ⇔ max(tAs, tBs) < min(tAe, tBe)"

$lines = @(
    "Search algorithm:",
    "Temporal DFS, temporal feasibility check based on constraint tightening:",
    $mathLine,
    $mathLine,
    "Join strategy: Nested-Loop;",
    "GraphNode-TSNode",
    "Each time a traversal goes through a virtual edge"
)

$tr.Text = [string]::Join("`r", $lines)

# Re-apply the outline (indent) levels that differ from the default (level 1 / IndentLevel=1):
$tr.Paragraphs(2, 1).IndentLevel = 2   # Temporal DFS, ...
$tr.Paragraphs(3, 1).IndentLevel = 2   # [tAs, tAe) ...
$tr.Paragraphs(4, 1).IndentLevel = 2   # [tAs, tAe) ... (duplicate)
$tr.Paragraphs(7, 1).IndentLevel = 2   # Each time a traversal goes through a virtual edge

# Re-apply bold on the one bullet that keeps it:
$tr.Paragraphs(6, 1).Font.Bold = $true  # GraphNode-TSNode
